$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing header cell (H1) onto new header cells I1, J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
